$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) AMR sheet comment text update (cell J1)
# ---------------------------------------------------------------------------
$wsAMR = $wb.Worksheets.Item("AMR")
$cmt = $wsAMR.Range("J1").Comment
$cmt.Text("Justin Replogle:`nHealth expenditure is not dependent on AMR scenario. NOTE Joao's spreadsheet shows these as positive numbers. I have made them negative and reversed the 5% and 95% levels to match the production losses.`n")

# ---------------------------------------------------------------------------
# 2) AMR sheet: swap K/L columns (health expenditure 5%/95% levels reversed)
#    for rows 6, 10 and 14, restyle J:L on those rows, and add matching
#    blank formatted cells on rows 7-9 and 11-13.
# ---------------------------------------------------------------------------

# Capture the DKK (no-decimals) cell style from J6 before changing it, so we
# can re-apply the plain DKK-with-decimals style afterwards.
$plainDkkStyle = $wsAMR.Range("M6").Style

foreach ($row in 6, 10, 14) {
    $kCell = $wsAMR.Range("K$row")
    $lCell = $wsAMR.Range("L$row")
    $kVal = $kCell.Value2
    $lVal = $lCell.Value2
    $kCell.Value2 = $lVal
    $lCell.Value2 = $kVal

    $wsAMR.Range("J$row:L$row").Style = $plainDkkStyle
}

# Re-assert the summary formulas explicitly (values recompute automatically,
# this just matches the literal (non-shared) formula text in the target).
$wsAMR.Range("I6").Formula = "=SUM(I3:I5)"
$wsAMR.Range("N6").Formula = "=H6+K6"
$wsAMR.Range("O6").Formula = "=I6+L6"

$blankStyle = $wsAMR.Range("D7").Style
foreach ($row in 7, 8, 9, 11, 12, 13) {
    $rng = $wsAMR.Range("J$row:L$row")
    $rng.Style = $blankStyle
    $rng.ClearContents()
}

# ---------------------------------------------------------------------------
# 3) Selections / active sheet state
# ---------------------------------------------------------------------------
$wsFarm = $wb.Worksheets.Item("Farm summary")
$wsFarm.Range("G16").Select()

$wsAHLE = $wb.Worksheets.Item("AHLE")
$wsAHLE.Range("C43").Select()

$wsAMR.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$wsAMR.Range("L15").Select()
